$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D (Price, stored as text) and E (Volume(1h)) values.
# Row 14's E value is unchanged in the source diff, so it is omitted here.
$updates = @{
  2  = @{ D = "22.432.34";   E = "  +9.28%  " }
  3  = @{ D = "1.603.38";    E = "  +8.75%  " }
  4  = @{ D = "1.006";       E = "  -0.27%  " }
  5  = @{ D = "306.63";      E = "  +10.41%  " }
  6  = @{ D = "0.9935";      E = "  +4.01%  " }
  7  = @{ D = "0.3673";      E = "  +1.59%  " }
  8  = @{ D = "0.3389";      E = "  +10.96%  " }
  9  = @{ D = "42.32";       E = "  +7.34%  " }
  10 = @{ D = "1.133";       E = "  +7.19%  " }
  11 = @{ D = "0.07042";     E = "  +6.01%  " }
  12 = @{ D = "1.003";       E = "  +0.03%  " }
  13 = @{ D = "19.71";       E = "  +8.83%  " }
  14 = @{ D = "5.909" }
  15 = @{ D = "6.621";       E = "  +7.10%  " }
  16 = @{ D = "1.605.32";    E = "  +8.76%  " }
  17 = @{ D = "0.00001081";  E = "  +5.04%  " }
  18 = @{ D = "0.9937";      E = "  +4.01%  " }
  19 = @{ D = "0.06625";     E = "  +11.32%  " }
  20 = @{ D = "77.91";       E = "  +12.73%  " }
  21 = @{ D = "5.999";       E = "  +9.21%  " }
  22 = @{ D = "16.00";       E = "  +10.48%  " }
  23 = @{ D = "11.83";       E = "  +6.15%  " }
  24 = @{ D = "22.489.50";   E = "  +9.35%  " }
  25 = @{ D = "2.395";       E = "  +5.89%  " }
  26 = @{ D = "2.569";       E = "  +20.69%  " }
  27 = @{ D = "149.45";      E = "  +4.53%  " }
  28 = @{ D = "19.54";       E = "  +13.88%  " }
  29 = @{ D = "1.788.76";    E = "  +9.37%  " }
  30 = @{ D = "122.89";      E = "  +8.13%  " }
  31 = @{ B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "4.036";  E = "  +2.46%  " }
  32 = @{ B = "Filecoin";   C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil";      D = "6.131";  E = "  +22.36%  " }
  33 = @{ D = "0.9496";      E = "  +17.56%  " }
  34 = @{ D = "1.694";       E = "  +12.22%  " }
  35 = @{ D = "0.08227";     E = "  +3.01%  " }
  36 = @{ D = "11.96";       E = "  +15.59%  " }
  37 = @{ D = "5.217";       E = "  +10.42%  " }
  38 = @{ D = "1.270";       E = "  +3.64%  " }
  39 = @{ D = "8.593";       E = "  +15.27%  " }
  40 = @{ D = "0.06078";     E = "  +3.77%  " }
  41 = @{ D = "0.02210";     E = "  +7.91%  " }
  42 = @{ D = "0.2020";      E = "  +7.61%  " }
  43 = @{ D = "0.9931";      E = "  +3.85%  " }
  44 = @{ D = "0.5902";      E = "  +11.56%  " }
  45 = @{ D = "3.848";       E = "  +9.23%  " }
  46 = @{ D = "13.11";       E = "  +7.52%  " }
  47 = @{ D = "0.5685";      E = "  +9.51%  " }
  48 = @{ D = "127.09";      E = "  +7.55%  " }
  49 = @{ D = "1.959";       E = "  +8.01%  " }
  50 = @{ D = "0.06815";     E = "  +5.28%  " }
  51 = @{ D = "73.42";       E = "  +8.89%  " }
}

foreach ($row in $updates.Keys) {
  $vals = $updates[$row]

  if ($vals.ContainsKey("B")) {
    $ws.Cells.Item($row, 2).Value = $vals["B"]
  }
  if ($vals.ContainsKey("C")) {
    $ws.Cells.Item($row, 3).Value = $vals["C"]
  }
  if ($vals.ContainsKey("D")) {
    # Price column values look numeric (e.g. "16.00", "1.270") but must stay
    # exact text, so force Text format before assigning, then restore the
    # default "Normal" cell style so no spurious style index is introduced.
    $dcell = $ws.Cells.Item($row, 4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $vals["D"]
    $dcell.Style = "Normal"
  }
  if ($vals.ContainsKey("E")) {
    $ws.Cells.Item($row, 5).Value = $vals["E"]
  }
}
